# ---------------------------------------------------------------------------
# baddata2015test.xlsx edit
#
# Summary of the change (see commit message / xml diff):
#   - Split-time number format is consolidated from two builtin time formats
#     (h:mm / [h]:mm:ss) into a single custom format "[h]:mm:ss;@" applied to
#     every timing cell (H3:O15).
#   - A number of split times in columns H, J, K, L (and a few M/N/O cells)
#     are corrected/recomputed.
#   - Two bib labels (row 7 and row 13, columns A/B) are renamed.
#   - A handful of cells that used to hold a (wrong) carried-over value are
#     cleared, and a few new, still-empty, but formatted cells appear at the
#     edge of the used range (O12:O14, K15:O15) because the format now
#     covers the full H3:O15 block.
#   - The active selection moves from A16 to A8.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the number format for every split-time cell -----------------
# Replaces the old mixed h:mm / [h]:mm:ss builtin formats with a single
# custom format, and (as a side effect) materializes formatted-but-empty
# cells at the bottom/right edge of the data block (O12:O14, K15:O15).
$ws.Range("H3:O15").NumberFormat = "[h]:mm:ss;@"

# --- Corrected split times (columns H/J/K/L) --------------------------------
$ws.Range("J3").Value = 0.4458333333333333
$ws.Range("K3").Value = 0.4486111111111111
$ws.Range("L3").Value = 0.6944444444444445

$ws.Range("J4").Value = 0.5298611111111111
$ws.Range("K4").Value = 0.53125
$ws.Range("L4").Value = 0.7013888888888888

$ws.Range("L5").Value = 0.7944444444444444

$ws.Range("J6").Value = 0.5652777777777778
$ws.Range("K6").Value = 0.56875
$ws.Range("L6").Value = 0.8263888888888888

$ws.Range("H7").Value = 0.020833333333333332
$ws.Range("J7").Value = 1.1354166666666667
$ws.Range("K7").Value = 1.1458333333333333
$ws.Range("L7").Value = 1.4590277777777778

$ws.Range("J8").Value = 1.1458333333333333
$ws.Range("K8").Value = 1.15625
$ws.Range("L8").Value = 1.4722222222222223

$ws.Range("J9").Value = 1.20625
$ws.Range("K9").Value = 1.215277777777778
$ws.Range("L9").Value = 1.5569444444444445

$ws.Range("J10").Value = 1.2249999999999999
$ws.Range("L10").Value = 1.3979166666666665

$ws.Range("J11").Value = 1.3166666666666667
$ws.Range("K11").Value = 1.3222222222222222
$ws.Range("L11").Value = 1.5972222222222223

$ws.Range("J12").Value = 1.2972222222222223
$ws.Range("K12").Value = 1.3055555555555556
$ws.Range("L12").Value = 1.5722222222222222

$ws.Range("H13").Value = 0.041666666666666664
$ws.Range("J13").Value = 1.3833333333333335
$ws.Range("K13").Value = 1.3875
$ws.Range("L13").Value = 1.6888888888888889

$ws.Range("J14").Value = 1.3375000000000001
$ws.Range("K14").Value = 1.3395833333333333
$ws.Range("L14").Value = 1.715972222222222

# --- Cells that no longer have (bogus) values; formatting is retained ------
$ws.Range("K10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("O10").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# --- Relabel rows 7 & 13 (bib name fields) ----------------------------------
$ws.Range("B7").Value = "Offset"
$ws.Range("A13").Value = "Hour"
$ws.Range("A7").Value = "Half-Hour"
$ws.Range("B13").Value = "Offset"

# --- Move the active selection from A16 to A8 -------------------------------
$ws.Range("A8").Select()
